$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Insert a new row at position 23 for "KATIVAROX SYRUP" (alphabetically
#    sorted between FLOPADEX and KETOLAC). This shifts rows 23-45 down to
#    24-46, preserving all their content, formatting and merges.
# ------------------------------------------------------------------
$ws.Rows.Item(23).Insert()

# Copy formatting + merges from the row right below (the shifted-down
# former row 23) into the newly-created blank row 23.
$ws.Rows.Item(24).Copy($ws.Rows.Item(23))
$ws.Rows.Item(23).RowHeight = 24.75

# Fill in the new row's data (text-like numeric values rely on the
# existing cell style's Text number format (numFmtId 49) to be stored
# as shared-string text rather than numbers).
$ws.Cells.Item(23,1).Value = 17
$ws.Cells.Item(23,3).Value = "KATIVAROX SYRUP"
$ws.Cells.Item(23,8).Value = "0:0"
$ws.Cells.Item(23,12).Value = "0"
$ws.Cells.Item(23,14).Value = "122.00"
$ws.Cells.Item(23,16).Value = "122.0000"
$ws.Cells.Item(23,17).Value = "1:0"

# ------------------------------------------------------------------
# 2. Renumber the "م" (row index) column for every shifted product row
#    (now rows 24-44): each one needs to be incremented by 1.
# ------------------------------------------------------------------
for ($r = 24; $r -le 44; $r++) {
    $ws.Cells.Item($r,1).Value = $ws.Cells.Item($r,1).Value2 + 1
}

# ------------------------------------------------------------------
# 3. Update the VIDROP row (now row 37) with its new stats.
# ------------------------------------------------------------------
$ws.Cells.Item(37,8).Value = "3:0"
$ws.Cells.Item(37,16).Value = "78.0000"
$ws.Cells.Item(37,17).Value = "3:0"

# ------------------------------------------------------------------
# 4. Update the totals row (now row 45) with the new summed value.
# ------------------------------------------------------------------
$ws.Cells.Item(45,16).Value = 2429.5749999999998

# ------------------------------------------------------------------
# 5. Update the footer timestamp (now row 46) from 7:56 PM to 7:57 PM.
# ------------------------------------------------------------------
$ws.Cells.Item(46,1).Value = "Friday, 25 July, 2025 7:57 PM"
